# RPA datasets push 2024-05-18
# Inserts two new IPO-subscription rows (KB제28호스팩, 아이씨티케이) at the
# top of the data table (rows 2-3), pushing the existing rows down by two
# (old row 2 -> row 4, ... old row 14 -> row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2..14) down to (4..16), keeping their
# original formatting intact.
$ws.Rows("2:3").Insert()

# The two blank rows created by Insert() inherit the header row's
# formatting (bold/border/center). Strip that back to the plain/default
# style used by every other data row.
$ws.Range("A2:T3").ClearFormats()

# --- Row 2: KB제28호스팩 ---------------------------------------------
$ws.Cells.Item(2, 1).Value  = "'2024-05-07"
$ws.Cells.Item(2, 2).Value  = "KB제28호스팩"
$ws.Cells.Item(2, 3).Value  = "KB"
$ws.Cells.Item(2, 4).Value  = "'2024-05-10"
$ws.Cells.Item(2, 5).Value  = "'2024-05-17"
$ws.Cells.Item(2, 6).Value  = 10000000
$ws.Cells.Item(2, 7).Value  = 5000000
$ws.Cells.Item(2, 8).Value  = "-"
$ws.Cells.Item(2, 9).Value  = 2000
$ws.Cells.Item(2, 10).Value = 2000
$ws.Cells.Item(2, 11).Value = "-"
$ws.Cells.Item(2, 12).Value = 2000
$ws.Cells.Item(2, 13).Value = "-"
$ws.Cells.Item(2, 14).Value = "-"
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = "-"
$ws.Cells.Item(2, 17).Value = "-"
$ws.Cells.Item(2, 18).Value = "713 : 1"
$ws.Cells.Item(2, 19).Value = "-"
$ws.Cells.Item(2, 20).Value = "-"

# --- Row 3: 아이씨티케이 ---------------------------------------------
$ws.Cells.Item(3, 1).Value  = "'2024-05-07"
$ws.Cells.Item(3, 2).Value  = "아이씨티케이"
$ws.Cells.Item(3, 3).Value  = "NH"
$ws.Cells.Item(3, 4).Value  = "'2024-05-10"
$ws.Cells.Item(3, 5).Value  = "'2024-05-17"
$ws.Cells.Item(3, 6).Value  = 39400000
$ws.Cells.Item(3, 7).Value  = 1970000
$ws.Cells.Item(3, 8).Value  = "-"
$ws.Cells.Item(3, 9).Value  = 13000
$ws.Cells.Item(3, 10).Value = 16000
$ws.Cells.Item(3, 11).Value = "-"
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = "-"
$ws.Cells.Item(3, 14).Value = "-"
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = "-"
$ws.Cells.Item(3, 17).Value = "-"
$ws.Cells.Item(3, 18).Value = "1108 : 1"
$ws.Cells.Item(3, 19).Value = "-"
$ws.Cells.Item(3, 20).Value = "-"

# The apostrophe-prefixed date strings above were entered as text (to stop
# Excel from auto-converting them to date serials) which leaves a
# "quote prefix" style on those cells - reset rows 2:3 back to the plain
# "Normal" style so they match every other (unstyled) data row.
$ws.Range("A2:T3").Style = "Normal"

Write-Output "done"
